# Update countries & provincias Spain
# Daily COVID-19 data refresh: updated case counters for several countries
# and re-sorted "Trinidad yTobago" to right after "Republica del Chad".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Estados Unidos (row 4) ---
$ws.Cells.Item(4, 2).Value = 5809027
$ws.Cells.Item(4, 3).Value = 12300
$ws.Cells.Item(4, 4).Value = 3127648
$ws.Cells.Item(4, 5).Value = 2501893
$ws.Cells.Item(4, 7).Value = 286
$ws.Cells.Item(4, 8).Value = 179486

# --- India (row 6) ---
$ws.Cells.Item(6, 2).Value = 3038013
$ws.Cells.Item(6, 3).Value = 64645
$ws.Cells.Item(6, 4).Value = 2273973
$ws.Cells.Item(6, 5).Value = 707248
$ws.Cells.Item(6, 7).Value = 864
$ws.Cells.Item(6, 8).Value = 56792

# --- Chile (row 13) ---
$ws.Cells.Item(13, 2).Value = 395708
$ws.Cells.Item(13, 3).Value = 1939
$ws.Cells.Item(13, 4).Value = 369730
$ws.Cells.Item(13, 5).Value = 15186
$ws.Cells.Item(13, 7).Value = 69
$ws.Cells.Item(13, 8).Value = 10792

# --- Argentina (row 15) ---
$ws.Cells.Item(15, 4).Value = 245781
$ws.Cells.Item(15, 5).Value = 76467
$ws.Cells.Item(15, 7).Value = 65
$ws.Cells.Item(15, 8).Value = 6795

# --- Reino Unido (row 16) ---
$ws.Cells.Item(16, 2).Value = 324601
$ws.Cells.Item(16, 3).Value = 1288
$ws.Cells.Item(16, 7).Value = 18
$ws.Cells.Item(16, 8).Value = 41423

# --- Italia (row 20) ---
$ws.Cells.Item(20, 2).Value = 258136
$ws.Cells.Item(20, 3).Value = 1071
$ws.Cells.Item(20, 4).Value = 205203
$ws.Cells.Item(20, 5).Value = 17503
$ws.Cells.Item(20, 7).Value = 3
$ws.Cells.Item(20, 8).Value = 35430

# --- Canada (row 27) ---
$ws.Cells.Item(27, 2).Value = 124481
$ws.Cells.Item(27, 3).Value = 109
$ws.Cells.Item(27, 4).Value = 110738
$ws.Cells.Item(27, 5).Value = 4678
$ws.Cells.Item(27, 7).Value = 1
$ws.Cells.Item(27, 8).Value = 9065

# --- Republica Dominicana (row 35) ---
$ws.Cells.Item(35, 2).Value = 90561
$ws.Cells.Item(35, 3).Value = 694
$ws.Cells.Item(35, 4).Value = 59949
$ws.Cells.Item(35, 5).Value = 29058
$ws.Cells.Item(35, 7).Value = 21
$ws.Cells.Item(35, 8).Value = 1554

# --- Guatemala (row 44) ---
$ws.Cells.Item(44, 2).Value = 67856
$ws.Cells.Item(44, 3).Value = 915
$ws.Cells.Item(44, 4).Value = 56277
$ws.Cells.Item(44, 5).Value = 8999
$ws.Cells.Item(44, 7).Value = 48
$ws.Cells.Item(44, 8).Value = 2580

# --- Azerbaiyan (row 64) ---
$ws.Cells.Item(64, 2).Value = 35105
$ws.Cells.Item(64, 3).Value = 184
$ws.Cells.Item(64, 4).Value = 32842
$ws.Cells.Item(64, 5).Value = 1748
$ws.Cells.Item(64, 7).Value = 3
$ws.Cells.Item(64, 8).Value = 515

# --- Madagascar (row 82) ---
$ws.Cells.Item(82, 2).Value = 14277
$ws.Cells.Item(82, 3).Value = 59
$ws.Cells.Item(82, 4).Value = 13332
$ws.Cells.Item(82, 5).Value = 767

# --- Haiti (row 99) ---
$ws.Cells.Item(99, 2).Value = 8050
$ws.Cells.Item(99, 3).Value = 34
$ws.Cells.Item(99, 5).Value = 2407

# --- Re-sort "Trinidad yTobago" up to right after "Republica del Chad" ---
# Before: 162 Santo Tome y Principe, 163 Guyana, 164 Trinidad yTobago
# After:  162 Trinidad yTobago (new data), 163 Santo Tome y Principe (old data), 164 Guyana (old data)

$ws.Cells.Item(163, 1).Value = "Santo Tome y Principe"
$ws.Cells.Item(163, 2).Value = 891
$ws.Cells.Item(163, 3).Value = 0
$ws.Cells.Item(163, 4).Value = 830
$ws.Cells.Item(163, 5).Value = 46
$ws.Cells.Item(163, 6).Value = 0
$ws.Cells.Item(163, 7).Value = 0
$ws.Cells.Item(163, 8).Value = 15

$ws.Cells.Item(164, 1).Value = "Guyana"
$ws.Cells.Item(164, 2).Value = 881
$ws.Cells.Item(164, 3).Value = 0
$ws.Cells.Item(164, 4).Value = 433
$ws.Cells.Item(164, 5).Value = 418
$ws.Cells.Item(164, 6).Value = 0
$ws.Cells.Item(164, 7).Value = 0
$ws.Cells.Item(164, 8).Value = 30

$ws.Cells.Item(162, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(162, 2).Value = 899
$ws.Cells.Item(162, 3).Value = 35
$ws.Cells.Item(162, 4).Value = 165
$ws.Cells.Item(162, 5).Value = 721
$ws.Cells.Item(162, 6).Value = 0
$ws.Cells.Item(162, 7).Value = 1
$ws.Cells.Item(162, 8).Value = 13

# --- Burundi (row 172) ---
$ws.Cells.Item(172, 2).Value = 429
$ws.Cells.Item(172, 3).Value = 3
$ws.Cells.Item(172, 5).Value = 92
